$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2: drop three addresses from the composed route text and
#     recompute dependent numeric columns ---
$ws.Range("D2").Value = "П-834 / Нижний Новгород г., 1-й Кемеровский пер.2; П-13 / Нижний Новгород г., 2-й Осташковский пер.1; П-325 / Нижний Новгород г., 30 лет Октября ул.2; П-331 / Нижний Новгород г., 40 лет Октября ул.26/1; П-696 / Нижний Новгород г., 40 лет Победы ул.1; "
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 22
$ws.Range("J2").Value = 2250.6
$ws.Range("K2").Value = 24.2
$ws.Range("M2").Value = 19.40176032355866
$ws.Range("Q2").Value = 0.3666666666666666
$ws.Range("S2").Value = 2.290857459119211
$ws.Range("T2").Value = 3.090857459119212
$ws.Range("U2").Value = 48
$ws.Range("V2").Value = 85.44023247062256
$ws.Range("W2").Value = 185.4514475471527

# --- Insert three new rows before the old row 4 (which carried the КГМ
#     aggregate route). These new rows represent the three stops that were
#     pulled out of the aggregate route in row 2 and the old row 4 route,
#     each becoming its own single-stop route. The old row 4 content shifts
#     down to row 7 automatically. ---
$ws.Rows("4:6").Insert()

# New row 4
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "приокский район"
$ws.Range("D4").Value = "П-709 / Нижний Новгород г., 40 лет Октября ул.15 к2"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = "КАМАЗ 43255-3010-69, МК-4512-04"
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 744
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 34.38357677950807
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 34.38357677950807
$ws.Range("O4").Value = 40
$ws.Range("P4").Value = 15
$ws.Range("Q4").Value = 0.25
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 1.146119225983603
$ws.Range("T4").Value = 1.646119225983603
$ws.Range("U4").Value = 30
$ws.Range("V4").Value = 68.76715355901614
$ws.Range("W4").Value = 98.76715355901615

# New row 5
$ws.Range("A5").Value = ""
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "приокский район"
$ws.Range("D5").Value = "П-326 / Нижний Новгород г., 40 лет Октября ул.15к1"
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = "КАМАЗ 43255-3010-69, МК-4512-04"
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 744
$ws.Range("K5").Value = 8
$ws.Range("L5").Value = 34.38357677950807
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 34.38357677950807
$ws.Range("O5").Value = 40
$ws.Range("P5").Value = 15
$ws.Range("Q5").Value = 0.25
$ws.Range("R5").Value = 0.25
$ws.Range("S5").Value = 1.146119225983603
$ws.Range("T5").Value = 1.646119225983603
$ws.Range("U5").Value = 30
$ws.Range("V5").Value = 68.76715355901614
$ws.Range("W5").Value = 98.76715355901615

# New row 6
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "приокский район"
$ws.Range("D6").Value = "П-1260 / Нижний Новгород г., 40 лет Октября ул.7Б"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = "КАМАЗ 43255-3010-69, МК-4512-04"
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 744
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 34.38357677950807
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 34.38357677950807
$ws.Range("O6").Value = 40
$ws.Range("P6").Value = 15
$ws.Range("Q6").Value = 0.25
$ws.Range("R6").Value = 0.25
$ws.Range("S6").Value = 1.146119225983603
$ws.Range("T6").Value = 1.646119225983603
$ws.Range("U6").Value = 30
$ws.Range("V6").Value = 68.76715355901614
$ws.Range("W6").Value = 98.76715355901615
